# The workbook originally has three columns:
#   A = Code (crime code), B = Name (powiat/county name), C = Crimes (count)
# The commit removes the "Name" column (B) entirely - the county-name text
# data (and its now-unused shared strings) goes away, and the former
# "Crimes" count column (C) shifts left to become the new column B.
# Net result: two columns remain - A = Code, B = Crimes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Name" column; Excel shifts column C (Crimes counts)
# left into B automatically, carrying its values/number formatting with it.
$ws.Columns("B").Delete()

# Leave the selection on the (now single) data column, matching the
# post-edit UI state (B1:B1048576, i.e. the whole column selected).
$ws.Range("B1:B1048576").Select() | Out-Null
